$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.166.00"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.652.95"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'217.74"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'0.5301"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.06325"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "'0.07806"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "1.660.34"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "1.881.14"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "'0.5490"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "0.0₅8180"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "'65.37"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "26.148.64"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'4.593"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "'191.06"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'6.004"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'145.27"
$ws.Range("E25").Value = "  +3.98%  "
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "'15.98"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'1.471"
$ws.Range("E29").Value = "  +3.68%  "
$ws.Range("D30").Value = "'0.05740"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'3.548"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "'3.263"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "'1.590"
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("D35").Value = "'2.803"
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").Value = "'2.423"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "'0.9488"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "'0.5731"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.795"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8493"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "1.039.54"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").Value = "'103.84"
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'56.69"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").Value = "'1.006"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4356"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").Value = "'0.05154"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'7.846"
